# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
#
# The match rows in this sheet got re-ordered upstream (same matches, rows
# shuffled). Columns A (row-sequence id), C (Div) and D (Date) stay put;
# everything else (B id .. AD PL_AhUnder) needs to move between rows
# following the cycles below. We snapshot every source row's values BEFORE
# writing anything, so rows inside the same cycle don't clobber each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-match data (everything except A/C/D).
$cols = @(2, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30)

# Each cycle lists rows in "receives-from" order: row[i] gets the data that
# currently lives in row[i+1] (wrapping around), i.e. target <- source.
$cycles = @(
    @(15, 17, 16),
    @(81, 83),
    @(129, 133),
    @(131, 132),
    @(154, 155)
)

foreach ($cycle in $cycles) {
    # Snapshot current values for every row in this cycle first.
    $snapshots = @{}
    foreach ($r in $cycle) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
        }
        $snapshots[$r] = $rowVals
    }

    # row[i] <- row[i+1] (the row after it in the list, wrapping to the front)
    $n = $cycle.Count
    for ($i = 0; $i -lt $n; $i++) {
        $target = $cycle[$i]
        $source = $cycle[($i + 1) % $n]
        $srcVals = $snapshots[$source]
        foreach ($c in $cols) {
            $ws.Cells.Item($target, $c).Value2 = $srcVals[$c]
        }
    }
}
